$d = $word.ActiveDocument

# The two Pearson logo pictures (in the footers) and the two BTec logo
# pictures (in the headers) need their inline-picture "name" (docPr/cNvPr
# @name) swapped: the Pearson logos go from "image2.png" to "image1.png",
# and the BTec logos go from "image1.jpg" to "image2.jpg". The displayed
# description (@descr) and the underlying media parts / relationships are
# left untouched.
#
# InlineShape.Name is not wired up for round-tripping in this host, so we
# go through the document's flat WordOpenXML (which does include the
# header/footer parts) and patch the attribute text directly, then write
# the whole thing back.

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml
